$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.960.76"
$ws.Range("E2").Value = "  +0.53%  "

$ws.Range("D3").Value = "3.144.42"
$ws.Range("E3").Value = "  -0.98%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.32"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.73"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.32%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.576"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -5.06%  "

$ws.Range("D9").Value = "3.158.33"
$ws.Range("E9").Value = "  -0.85%  "

$ws.Range("E10").Value = "  -1.26%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.64"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.77%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.383"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.15%  "

$ws.Range("D13").Value = "3.694.11"
$ws.Range("E13").Value = "  -0.89%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.126"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.09%  "

$ws.Range("D15").Value = "64.961.07"
$ws.Range("E15").Value = "  +0.50%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.16"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.38%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.161.86"
$ws.Range("E17").Value = "  -0.64%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000156"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.21%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "411.65"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.59%  "

$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.25"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.69%  "

$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.56"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.04%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.08"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.15%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.16%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.86"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.81%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.485"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.80%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.194"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -4.77%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000104"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.53%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.22"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +3.98%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.994"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.27%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.81"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.07%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.09%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.36"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.93%  "

$ws.Range("E33").Value = "  -2.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "163.77"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +4.92%  "

$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.15"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.26%  "

$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.28"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.38%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.37"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.35%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.69"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.60%  "

$ws.Range("D39").Value = "2.605.19"
$ws.Range("E39").Value = "  -3.77%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.98"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.30%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.16"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.86%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.34"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.696"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.93%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0621"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.20%  "

$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.51"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.65%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0258"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.49%  "

$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.28"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -5.34%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "290.07"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.55%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.995"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.26%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0976"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.74%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.92"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.99%  "
